# Generate Report for Handoff
#
# The localization-status report marks the row for
# "a81f6c4a-9f7b-4245-8ab5-b0948593b50b.md" as handed off:
#   - Overview sheet: Status columns (zh-cn / de-de) go from
#     "In Translation" to "Ready for handoff".
#   - zh-cn sheet: Status goes to "Ready for handoff" and the
#     "Latest Handoff Datetime" is refreshed.
#   - de-de sheet: Status goes to "Ready for handoff" and the
#     "Latest Handoff Datetime" is refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet - row 7 corresponds to a81f6c4a-9f7b-4245-8ab5-b0948593b50b.md
$overview.Range("B7").Value = "Ready for handoff"
$overview.Range("C7").Value = "Ready for handoff"

# zh-cn sheet - row 7 corresponds to a81f6c4a-9f7b-4245-8ab5-b0948593b50b.md
$zhcn.Range("B7").Value = "Ready for handoff"
$zhcn.Range("D7").Value = "2016-03-08 20:48:58"

# de-de sheet - row 7 corresponds to a81f6c4a-9f7b-4245-8ab5-b0948593b50b.md
$dede.Range("B7").Value = "Ready for handoff"
$dede.Range("D7").Value = "2016-03-08 20:49:06"
